$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "код идентификатор PLU"

# C3 (Товар 2 row): numeric 50 -> text "fd"
$ws.Range("C3").Value = "fd"

# C4 (Товар 3 row): text "fd" -> numeric 35
$ws.Range("C4").Value = 35

# C5 (Товар 4 row): numeric 5 -> numeric 20
$ws.Range("C5").Value = 20

# C6 (Товар 5 row): text "20.00" -> numeric 5
$ws.Range("C6").Value = 5

# Row 1 height returns to default (auto) after text no longer wraps to 3 lines
$ws.Rows.Item(1).AutoFit()

# Column A widened manually (no longer best-fit)
$ws.Columns.Item(1).ColumnWidth = 36.7109375

# Update the active selection to mirror the authored view state
$ws.Range("C12").Select()
